$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("items")

# Header row: normalize C1/D1 fill to match B1 (yellow)
$ws.Range("C1:D1").Interior.Color = $ws.Range("B1").Interior.Color

# Row 2: replace content (IFB washing machine -> Samsung AC), strip hyperlink styling
$ws.Range("A2").Value = "samsung 1.5 ton 5star"
$ws.Range("B2").Value = "Samsung"
$ws.Range("C2").Value = "1.5 Ton"

$ws.Range("A2").Style = "Normal"
$ws.Range("A2").NumberFormat = "0"
$ws.Range("B2").Style = "Normal"

# Row 3: D3 color change
$ws.Range("D3").Value = "Pink"

$null = $ws.Range("D3").Select()
